$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.955.15"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "2.417.17"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'562.67"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").Value = "'142.68"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D11").Value = "'5.19"
$ws.Range("E11").Value = "  -4.06%  "

$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("D13").Value = "'25.84"
$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").Value = "'0.0000172"
$ws.Range("E14").Value = "  -2.09%  "

$ws.Range("D15").Value = "2.853.87"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").Value = "61.869.76"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "2.420.99"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").Value = "'11.30"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "'323.03"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("E20").Value = "  -1.62%  "

$ws.Range("D21").Value = "'6.81"
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "'66.65"
$ws.Range("E23").Value = "  +1.86%  "

$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").Value = "'8.69"
$ws.Range("E25").Value = "  -3.87%  "

$ws.Range("D26").Value = "'550.00"
$ws.Range("E26").Value = "  -6.60%  "

$ws.Range("D27").Value = "2.536.69"
$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("E30").Value = "  -1.63%  "

$ws.Range("E31").Value = "  -4.42%  "

$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("E33").Value = "  -0.76%  "

$ws.Range("E34").Value = "  -4.28%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").Value = "'4.73"
$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("E37").Value = "  -1.52%  "

$ws.Range("D38").Value = "'153.44"
$ws.Range("E38").Value = "  +2.02%  "

$ws.Range("D39").Value = "'5.41"
$ws.Range("E39").Value = "  -5.07%  "

$ws.Range("D40").Value = "'18.54"
$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").Value = "'0.992"
$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("D43").Value = "'146.59"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("D44").Value = "'2.22"
$ws.Range("E44").Value = "  -5.17%  "

$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("E46").Value = "  -2.73%  "

$ws.Range("E47").Value = "  -1.93%  "

$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("D51").Value = "'11.56"
$ws.Range("E51").Value = "  +0.69%  "
